$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measurements for row 7 (the box that was re-measured)
$ws.Range("C7").Value = 57.6
$ws.Range("D7").Value = 37.4
$ws.Range("E7").Value = 37.9

# Add label in F7 referencing the new shared string "Mekkeske"
$ws.Range("F7").Value = "Mekkeske"

# Add new row with value in C8
$ws.Range("C8").Value = 0.37

# Update selection to match the new active cell
$ws.Range("F8").Select()
